# Weekly fruit/vegetable price update: insert a new weekly record as the
# first data row (row 34), pushing all subsequent rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 34 - shifts rows 34:100 down to 35:101
# and extends the sheet dimension to A1:R101.
$ws.Rows("34").Insert()

# Populate the newly inserted row with the new weekly price record.
$ws.Range("A34").Value = 8
$ws.Range("B34").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C34").Value = 'Coquimbo'
$ws.Range("D34").Value = (Get-Date -Year 2021 -Month 11 -Day 19 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E34").Value = 4
$ws.Range("F34").Value = 100112001
$ws.Range("G34").Value = 'Berenjena'
$ws.Range("H34").Value = 'Sin especificar'
$ws.Range("I34").Value = 'Primera'
$ws.Range("J34").Value = 560
$ws.Range("K34").Value = 8000
$ws.Range("L34").Value = 8500
$ws.Range("M34").Value = 8250
$ws.Range("N34").Value = '$/caja 60 unidades'
$ws.Range("O34").Value = 'Región de Arica y Parinacota'
$ws.Range("P34").Value = 138
$ws.Range("Q34").Value = 60
$ws.Range("R34").Value = 'Hortaliza'
